# Update query dashboard dan tampilan status evaluasi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Masa_Studi" column (J). This shifts HP_Ortu/HP_Mahasiswa/Email/Status/Evaluasi
# left by one (K->J, L->K, M->L, N->M, O->N).
$ws.Columns("J").Delete() | Out-Null

# Remove the "Evaluasi" column, which is now column N after the previous delete.
# This shifts Status (M) to stay at M (no further shift needed beyond removing N).
$ws.Columns("N").Delete() | Out-Null

# Update Total_SKS values that changed for a handful of students.
$ws.Range("I13").Value = 53
$ws.Range("I14").Value = 40
$ws.Range("I26").Value = 44
$ws.Range("I27").Value = 31

# Update the view/selection state to reflect the new scroll position and active cell.
$ws.Range("H7").Select() | Out-Null
